$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.747.06'
$ws.Range("E2").Value = '  -6.33%  '

# Row 3
$ws.Range("D3").Value = '2.463.54'
$ws.Range("E3").Value = '  -8.79%  '

# Row 4
$ws.Range("E4").Value = '  +0.25%  '

# Row 5
$ws.Range("D5").Formula = '''469.10'
$ws.Range("E5").Value = '  -6.57%  '

# Row 6
$ws.Range("D6").Formula = '''133.63'
$ws.Range("E6").Value = '  -4.49%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Formula = '''0.498'
$ws.Range("E8").Value = '  -5.58%  '

# Row 9
$ws.Range("D9").Value = '2.463.09'
$ws.Range("E9").Value = '  -9.24%  '

# Row 10
$ws.Range("D10").Formula = '''0.0964'
$ws.Range("E10").Value = '  -7.86%  '

# Row 11
$ws.Range("D11").Formula = '''5.35'
$ws.Range("E11").Value = '  -11.48%  '

# Row 12
$ws.Range("D12").Formula = '''0.318'
$ws.Range("E12").Value = '  -8.43%  '

# Row 13
$ws.Range("E13").Value = '  -3.59%  '

# Row 14
$ws.Range("D14").Value = '2.890.47'
$ws.Range("E14").Value = '  -9.06%  '

# Row 15
$ws.Range("D15").Value = '54.905.58'
$ws.Range("E15").Value = '  -6.17%  '

# Row 16
$ws.Range("D16").Formula = '''0.0000136'
$ws.Range("E16").Value = '  +1.47%  '

# Row 17
$ws.Range("D17").Formula = '''19.85'
$ws.Range("E17").Value = '  -7.61%  '

# Row 18
$ws.Range("D18").Value = '2.450.89'
$ws.Range("E18").Value = '  -9.67%  '

# Row 19
$ws.Range("D19").Formula = '''4.27'
$ws.Range("E19").Value = '  -9.92%  '

# Row 20
$ws.Range("D20").Formula = '''314.87'
$ws.Range("E20").Value = '  -6.00%  '

# Row 21
$ws.Range("D21").Formula = '''9.67'
$ws.Range("E21").Value = '  -11.14%  '

# Row 22
$ws.Range("D22").Formula = '''1.00'
$ws.Range("E22").Value = '  +0.66%  '

# Row 23
$ws.Range("D23").Formula = '''5.68'
$ws.Range("E23").Value = '  +0.87%  '

# Row 24
$ws.Range("D24").Formula = '''5.43'
$ws.Range("E24").Value = '  -12.62%  '

# Row 25
$ws.Range("D25").Formula = '''57.03'
$ws.Range("E25").Value = '  -9.44%  '

# Row 26
$ws.Range("D26").Formula = '''1.00'
$ws.Range("E26").Value = '  +0.90%  '

# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Formula = '''0.159'
$ws.Range("E27").Value = '  -8.49%  '

# Row 28
$ws.Range("B28").Value = 'Polygon'
$ws.Range("C28").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D28").Formula = '''0.389'
$ws.Range("E28").Value = '  -8.34%  '

# Row 29
$ws.Range("D29").Value = '2.549.17'
$ws.Range("E29").Value = '  -9.96%  '

# Row 30
$ws.Range("D30").Formula = '''7.21'
$ws.Range("E30").Value = '  -3.06%  '

# Row 31
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("D32").Value = '0.0₃0729'
$ws.Range("E32").Value = '  -10.89%  '

# Row 33
$ws.Range("D33").Formula = '''147.64'
$ws.Range("E33").Value = '  -2.45%  '

# Row 34
$ws.Range("D34").Formula = '''17.94'
$ws.Range("E34").Value = '  -6.21%  '

# Row 35
$ws.Range("D35").Formula = '''1.44'
$ws.Range("E35").Value = '  -9.75%  '

# Row 36
$ws.Range("D36").Formula = '''5.08'
$ws.Range("E36").Value = '  -5.68%  '

# Row 37
$ws.Range("D37").Formula = '''3.63'
$ws.Range("E37").Value = '  -12.80%  '

# Row 38
$ws.Range("E38").Value = '  -5.01%  '

# Row 39
$ws.Range("D39").Formula = '''0.806'
$ws.Range("E39").Value = '  -13.28%  '

# Row 40
$ws.Range("D40").Formula = '''0.998'
$ws.Range("E40").Value = '  +0.12%  '

# Row 41
$ws.Range("D41").Formula = '''33.04'
$ws.Range("E41").Value = '  -7.17%  '

# Row 42
$ws.Range("D42").Formula = '''0.602'
$ws.Range("E42").Value = '  +1.27%  '

# Row 43
$ws.Range("D43").Formula = '''0.0530'
$ws.Range("E43").Value = '  -5.02%  '

# Row 44
$ws.Range("D44").Formula = '''3.29'
$ws.Range("E44").Value = '  -7.51%  '

# Row 45
$ws.Range("D45").Formula = '''1.25'
$ws.Range("E45").Value = '  -9.12%  '

# Row 46
$ws.Range("D46").Formula = '''10.10'
$ws.Range("E46").Value = '  -2.46%  '

# Row 47
$ws.Range("D47").Value = '1.948.67'
$ws.Range("E47").Value = '  -10.51%  '

# Row 48
$ws.Range("D48").Formula = '''0.0890'
$ws.Range("E48").Value = '  +0.33%  '

# Row 49
$ws.Range("D49").Formula = '''0.0219'
$ws.Range("E49").Value = '  -2.83%  '

# Row 50
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Formula = '''4.31'
$ws.Range("E50").Value = '  -7.10%  '

# Row 51
$ws.Range("B51").Value = 'Bittensor'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D51").Formula = '''236.24'
$ws.Range("E51").Value = '  +5.98%  '
